# Swap the roster rows for (Anthony Lamb (TW) <-> Moses Moody) and
# (Andrew Wiggins <-> JaMychal Green) while keeping the leading index
# column (A) fixed to its row position.
#
# This mirrors the authoritative diff: columns B..K (No., Player, Pos,
# Ht, Wt, Birth Date, Unnamed:6, Exp, College, bbref url) trade places
# between the two row pairs, but column A (the 0-based row counter)
# stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($ws, $row1, $row2) {
    # Scratch cell, well outside the table (A1:K17), used as temporary
    # holding space for a true cell-to-cell copy/paste swap.
    $scratch = $ws.Range("Z1")
    # Plain value columns - direct Value2 read/write swap is fine here;
    # COM's natural type inference matches what was already stored for
    # these (numbers stay numbers, plain text stays plain text).
    $cols = @("B", "C", "D", "E", "F", "G", "H", "J", "K")

    $vals1 = @{}
    $vals2 = @{}
    foreach ($col in $cols) {
        $vals1[$col] = $ws.Range("$col$row1").Value2
        $vals2[$col] = $ws.Range("$col$row2").Value2
    }

    foreach ($col in $cols) {
        $ws.Range("$col$row1").Value2 = $vals2[$col]
        $ws.Range("$col$row2").Value2 = $vals1[$col]
    }

    # Column I ("Exp") holds numeric-looking text ("1", "2", "8", ...)
    # alongside non-numeric entries ("R") elsewhere in the same column,
    # so the source workbook stores the whole column as shared-string
    # text rather than numbers. A plain Value2 assignment would let
    # COM's numeric auto-detection silently retype "1"/"2" as numbers,
    # so swap these two cells with a real copy/paste (via a scratch
    # cell) instead - that carries the original text cell type across
    # without disturbing styles.
    $ws.Range("I$row1").Copy($scratch)
    $ws.Range("I$row2").Copy($ws.Range("I$row1"))
    $scratch.Copy($ws.Range("I$row2"))
    $scratch.ClearContents()
}

Swap-RowData $ws 8 9
Swap-RowData $ws 11 12
